$d = $word.ActiveDocument

function Expand-FldSimpleField($instr) {
    $escaped = $instr -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
        '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText>' + $escaped + '</w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Collect the paragraphs that hold a simple field (w:fldSimple) along with
# the field's instruction text, then replace each paragraph's range content
# with the expanded begin/instrText/separate/end run sequence. We gather
# first (fields shift indices/objects once content is replaced), then apply.
$targets = New-Object System.Collections.ArrayList
foreach ($f in $d.Fields) {
    if ($f.Kind -eq 1) {
        [void]$targets.Add(@{ Start = $f.Code.Start; End = $f.Code.End; Instr = $f.Code.Text.Trim() })
    }
}

foreach ($t in $targets) {
    $range = $d.Range($t.Start, $t.End)
    $para = $range.Paragraphs(1)
    [void]$para.Range.InsertXML((Expand-FldSimpleField $t.Instr))
}
